# Update column G ("K") values for data rows 2-48 in the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(0,2,0,1,1,1,1,0,2,2,0,1,1,1,4,1,2,4,2,1,1,1,1,1,0,2,0,4,1,1,3,1,0,1,0,1,3,4,0,3,2,0,5,2,2,0,2)

$row = 2
foreach ($v in $newValues) {
    $ws.Cells.Item($row, 7).Value = $v
    $row++
}
